$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.4587784673220199
$ws.Range("D2").Value = 0.2194335328818973
$ws.Range("E2").Value = 0.1452770811831083
$ws.Range("F2").Value = 0.8257246146173074
$ws.Range("G2").Value = 0.3995463891372708
$ws.Range("H2").Value = 0.5241085935554253
$ws.Range("I2").Value = 1.054375338444487
$ws.Range("J2").Value = 0.139081600977164
$ws.Range("M2").Value = 1.688862129694087
$ws.Range("O2").Value = 1.791829452560165
$ws.Range("B3").Value = 0.4029174815801184
$ws.Range("D3").Value = 0.2231179463799133
$ws.Range("E3").Value = 0.1487658639803442
$ws.Range("F3").Value = 0.824837045401928
$ws.Range("G3").Value = 0.3904856432766479
$ws.Range("H3").Value = 0.5243847230650829
$ws.Range("I3").Value = 0.961619056993996
$ws.Range("J3").Value = 0.1421067658121089
$ws.Range("M3").Value = 1.491933859935628
$ws.Range("O3").Value = 1.772863791967353
$ws.Range("B4").Value = 0.368474994892722
$ws.Range("D4").Value = 0.2255144771000523
$ws.Range("E4").Value = 0.151039282539327
$ws.Range("F4").Value = 0.8249840607728487
$ws.Range("G4").Value = 0.3853062208217324
$ws.Range("H4").Value = 0.5248901381313829
$ws.Range("I4").Value = 0.9050041978951811
$ws.Range("J4").Value = 0.1441171771289262
$ws.Range("M4").Value = 1.370363517415143
$ws.Range("O4").Value = 1.76272372683303
$ws.Range("B5").Value = 0.3544040846944654
$ws.Range("D5").Value = 0.2265248070752808
$ws.Range("E5").Value = 0.1519986560217337
$ws.Range("F5").Value = 0.8252173736551569
$ws.Range("G5").Value = 0.3832915930706804
$ws.Range("H5").Value = 0.525180415192338
$ws.Range("I5").Value = 0.8820202403759936
$ws.Range("J5").Value = 0.1449747945099702
$ws.Range("M5").Value = 1.320662736533194
$ws.Range("O5").Value = 1.75896861192021
$ws.Range("B6").Value = 0.3520655116020919
$ws.Range("D6").Value = 0.2266946071930338
$ws.Range("E6").Value = 0.1521599463065444
$ws.Range("F6").Value = 0.8252665705231124
$ws.Range("G6").Value = 0.382962850670566
$ws.Range("H6").Value = 0.5252337038514412
$ws.Range("I6").Value = 0.8782090928121988
$ws.Range("J6").Value = 0.1451195150332261
$ws.Range("M6").Value = 1.312400451077778
$ws.Range("O6").Value = 1.758367809997253
$ws.Range("B7").Value = 0.3682853716031502
$ws.Range("D7").Value = 0.225527966272848
$ws.Range("E7").Value = 0.1510520877288117
$ws.Range("F7").Value = 0.8249865060132962
$ws.Range("G7").Value = 0.3852786627062841
$ws.Range("H7").Value = 0.5248937117052321
$ws.Range("I7").Value = 0.9046938728842946
$ws.Range("J7").Value = 0.1441285880802745
$ws.Range("M7").Value = 1.369693876004305
$ws.Range("O7").Value = 1.762671559244751
$ws.Range("B8").Value = 0.4395478606005554
$ws.Range("D8").Value = 0.2206760222655664
$ws.Range("E8").Value = 0.1464527176387866
$ws.Range("F8").Value = 0.8252745939522796
$ws.Range("G8").Value = 0.396342253326182
$ws.Range("H8").Value = 0.5241340068263582
$ws.Range("I8").Value = 1.022324052184814
$ws.Range("J8").Value = 0.140092872153879
$ws.Range("M8").Value = 1.621100232345185
$ws.Range("O8").Value = 1.78497690312301
$ws.Range("B9").Value = 0.5781265642575306
$ws.Range("D9").Value = 0.2122286397507729
$ws.Range("E9").Value = 0.1384784527218361
$ws.Range("F9").Value = 0.8313585810149107
$ws.Range("G9").Value = 0.4211086567763118
$ws.Range("H9").Value = 0.5253162531437567
$ws.Range("I9").Value = 1.255600870901191
$ws.Range("J9").Value = 0.1333970069726931
$ws.Range("M9").Value = 2.108712590691823
$ws.Range("O9").Value = 1.840725731586275
$ws.Range("B10").Value = 0.6792036155396488
$ws.Range("D10").Value = 0.2066751752465059
$ws.Range("E10").Value = 0.1332612873910346
$ws.Range("F10").Value = 0.8392342331979421
$ws.Range("G10").Value = 0.4412128567615099
$ws.Range("H10").Value = 0.527824524396749
$ws.Range("I10").Value = 1.428497202878788
$ws.Range("J10").Value = 0.129226540876715
$ws.Range("M10").Value = 2.463450805740081
$ws.Range("O10").Value = 1.889104857935308
$ws.Range("B11").Value = 0.7250218973471192
$ws.Range("D11").Value = 0.2042909735779954
$ws.Range("E11").Value = 0.1310281137320004
$ws.Range("F11").Value = 0.8435655465445535
$ws.Range("G11").Value = 0.4507810795431055
$ws.Range("H11").Value = 0.5293240174856635
$ws.Range("I11").Value = 1.507465272804438
$ws.Range("J11").Value = 0.1274933242362692
$ws.Range("M11").Value = 2.624024511360403
$ws.Range("O11").Value = 1.912746947153465
$ws.Range("B12").Value = 0.7423482154043199
$ws.Range("D12").Value = 0.2034086161151993
$ws.Range("E12").Value = 0.1302027061228639
$ws.Range("F12").Value = 0.8453140411720312
$ws.Range("G12").Value = 0.4544656905315776
$ws.Range("H12").Value = 0.5299435582518157
$ws.Range("I12").Value = 1.537412456396339
$ws.Range("J12").Value = 0.1268606965100574
$ws.Range("M12").Value = 2.684710488730843
$ws.Range("O12").Value = 1.921936202998296
$ws.Range("B13").Value = 0.7386177650551531
$ws.Range("D13").Value = 0.2035977355070955
$ws.Range("E13").Value = 0.1303795704253137
$ws.Range("F13").Value = 0.8449326433947562
$ws.Range("G13").Value = 0.4536694071630905
$ws.Range("H13").Value = 0.5298078263360537
$ws.Range("I13").Value = 1.530960880130294
$ws.Range("J13").Value = 0.1269958879192146
$ws.Range("M13").Value = 2.671646084265149
$ws.Range("O13").Value = 1.919946589382619
$ws.Range("B14").Value = 0.7264478305904163
$ws.Range("D14").Value = 0.2042179705773943
$ws.Range("E14").Value = 0.1309598006313717
$ws.Range("F14").Value = 0.8437072215368744
$ws.Range("G14").Value = 0.4510829823865379
$ws.Range("H14").Value = 0.5293739500636718
$ws.Range("I14").Value = 1.509928180344644
$ws.Range("J14").Value = 0.1274408017734672
$ws.Range("M14").Value = 2.629019612450577
$ws.Range("O14").Value = 1.913498204079531
$ws.Range("B15").Value = 0.7189902377319299
$ws.Range("D15").Value = 0.2046005522199383
$ws.Range("E15").Value = 0.1313178479345058
$ws.Range("F15").Value = 0.8429707415237999
$ws.Range("G15").Value = 0.4495067260321974
$ws.Range("H15").Value = 0.5291149281609364
$ws.Range("I15").Value = 1.497050680483312
$ws.Range("J15").Value = 0.1277164151207231
$ws.Range("M15").Value = 2.602893919735806
$ws.Range("O15").Value = 1.909579227393692
$ws.Range("B16").Value = 0.6762059550069921
$ws.Range("D16").Value = 0.2068338524601732
$ws.Range("E16").Value = 0.133410058619484
$ws.Range("F16").Value = 0.8389662932640363
$ws.Range("G16").Value = 0.4405961038272892
$ws.Range("H16").Value = 0.5277337566516991
$ws.Range("I16").Value = 1.423342666385054
$ws.Range("J16").Value = 0.1293431206498212
$ws.Range("M16").Value = 2.452940441451972
$ws.Range("O16").Value = 1.887592805414897
$ws.Range("B17").Value = 0.6499171049117649
$ws.Range("D17").Value = 0.2082403432064197
$ws.Range("E17").Value = 0.1347295227962082
$ws.Range("F17").Value = 0.8367019134886959
$ws.Range("G17").Value = 0.435238383621666
$ws.Range("H17").Value = 0.5269783844717892
$ws.Range("I17").Value = 1.378204955459836
$ws.Range("J17").Value = 0.1303831337505059
$ws.Range("M17").Value = 2.360740760309142
$ws.Range("O17").Value = 1.874524480772891
$ws.Range("B18").Value = 0.6347812322232755
$ws.Range("D18").Value = 0.2090626910916846
$ws.Range("E18").Value = 0.1355016287318431
$ws.Range("F18").Value = 0.835469927872353
$ws.Range("G18").Value = 0.4321965206532781
$ws.Range("H18").Value = 0.5265776450184347
$ws.Range("I18").Value = 1.352272849441903
$ws.Range("J18").Value = 0.1309967456906165
$ws.Range("M18").Value = 2.307635077417842
$ws.Range("O18").Value = 1.867161652543672
$ws.Range("B19").Value = 0.6296538947531474
$ws.Range("D19").Value = 0.2093434188782446
$ws.Range("E19").Value = 0.1357653122032207
$ws.Range("F19").Value = 0.835064874061203
$ws.Range("G19").Value = 0.4311734095953454
$ws.Range("H19").Value = 0.5264477496790079
$ws.Range("I19").Value = 1.343497878663356
$ws.Range("J19").Value = 0.1312071493493789
$ws.Range("M19").Value = 2.289641700547065
$ws.Range("O19").Value = 1.864695081842257
$ws.Range("B20").Value = 0.6527171787537611
$ws.Range("D20").Value = 0.2080892355575479
$ws.Range("E20").Value = 0.1345876982944105
$ws.Range("F20").Value = 0.8369356669432904
$ws.Range("G20").Value = 0.4358046035291636
$ws.Range("H20").Value = 0.5270553029142775
$ws.Range("I20").Value = 1.383006859720325
$ws.Range("J20").Value = 0.1302708253389753
$ws.Range("M20").Value = 2.370563358310875
$ws.Range("O20").Value = 1.875899705128063
$ws.Range("B21").Value = 0.7300230944838404
$ws.Range("D21").Value = 0.204035236031288
$ws.Range("E21").Value = 0.1307888227057123
$ws.Range("F21").Value = 0.8440642125538886
$ws.Range("G21").Value = 0.4518410091440757
$ws.Range("H21").Value = 0.5294999851656286
$ws.Range("I21").Value = 1.516104822756347
$ws.Range("J21").Value = 0.1273094754458555
$ws.Range("M21").Value = 2.641543331128304
$ws.Range("O21").Value = 1.915385820276555
$ws.Range("B22").Value = 0.7804063443278437
$ws.Range("D22").Value = 0.2015051577107805
$ws.Range("E22").Value = 0.1284240959802059
$ws.Range("F22").Value = 0.8493548109438791
$ws.Range("G22").Value = 0.462679473436026
$ws.Range("H22").Value = 0.5313992405141192
$ws.Range("I22").Value = 1.603345883954262
$ws.Range("J22").Value = 0.1255122866399319
$ws.Range("M22").Value = 2.817943866935337
$ws.Range("O22").Value = 1.942571540708741
$ws.Range("B23").Value = 0.7535289655561996
$ws.Range("D23").Value = 0.2028445609502203
$ws.Range("E23").Value = 0.1296753614042974
$ws.Range("F23").Value = 0.8464731001840988
$ws.Range("G23").Value = 0.4568618672442
$ws.Range("H23").Value = 0.5303579281152508
$ws.Range("I23").Value = 1.556761052234663
$ws.Range("J23").Value = 0.1264587898118208
$ws.Range("M23").Value = 2.723861305318991
$ws.Range("O23").Value = 1.927935307278375
$ws.Range("B24").Value = 0.6514513340375174
$ws.Range("D24").Value = 0.2081575085371199
$ws.Range("E24").Value = 0.1346517750164526
$ws.Range("F24").Value = 0.8368297695163847
$ws.Range("G24").Value = 0.4355484960554605
$ws.Range("H24").Value = 0.5270204236381346
$ws.Range("I24").Value = 1.380835862081938
$ws.Range("J24").Value = 0.130321551090681
$ws.Range("M24").Value = 2.366122869308867
$ws.Range("O24").Value = 1.875277498047012
$ws.Range("B25").Value = 0.5407649891818096
$ws.Range("D25").Value = 0.2143993949582494
$ws.Range("E25").Value = 0.140523364711048
$ws.Range("F25").Value = 0.829117586449712
$ws.Range("G25").Value = 0.4140761861655449
$ws.Range("H25").Value = 0.5247092557513469
$ws.Range("I25").Value = 1.192223029368051
$ws.Range("J25").Value = 0.1350773524862809
$ws.Range("M25").Value = 1.977401309736734
$ws.Range("O25").Value = 1.824348299289767
